$d = $word.ActiveDocument
$full = $d.Content
$xml = $full.WordOpenXML

# 1. Remove the stray <w:lang w:val="en-US"/> elements from the title run
#    properties (there are exactly two, in the document title paragraph).
$xml = $xml.Replace('<w:lang w:val="en-US"/>', '')

# 2. Swap the single table's inside borders (insideH/insideV) for outside
#    borders (top/left/bottom/right).
$oldBorders = '<w:tblBorders><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders>'
$newBorders = '<w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders>'
$xml = $xml.Replace($oldBorders, $newBorders)

# 3. Shrink that table's single grid column from 9026 to 9016 twips.
$xml = $xml.Replace('<w:gridCol w:w="9026"/>', '<w:gridCol w:w="9016"/>')

# 4. Merge the leading " " run with the following label run in each of the
#    eight proprietor-detail paragraphs (Nome / E-mail / Endereço / CPF /
#    Fone / Cidade / CEP / UF), producing a single run whose text is the
#    space plus the original label text.
$mergePairs = @(
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Nome</w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> Nome</w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">E-mail: </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> E-mail: </w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Endereço: </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> Endereço: </w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">CPF:  </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> CPF:  </w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Fone: </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> Fone: </w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Cidade: </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> Cidade: </w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">CEP: </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> CEP: </w:t></w:r>'),
    @('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">UF: </w:t></w:r>',
      '<w:r><w:t xml:space="preserve"> UF: </w:t></w:r>')
)

foreach ($pair in $mergePairs) {
    $xml = $xml.Replace($pair[0], $pair[1])
}

$full.InsertXML($xml)
Write-Output "edit applied"
